$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "CreateContact"
$ws.Range("A2").Select()
